# Insert a new weekly price record for "Femacal de La Calera" (Espinaca)
# at row 468, pushing the existing rows 468:506 down to 469:507.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 468 (shifts rows 468-506 down to 469-507)
$ws.Rows.Item(468).Insert()

# Populate the newly inserted row with the new record's data
$ws.Cells.Item(468, 1).Value = 3
$ws.Cells.Item(468, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(468, 3).Value = "Coquimbo"
$ws.Cells.Item(468, 4).Value = 45013
$ws.Cells.Item(468, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(468, 5).Value = 5
$ws.Cells.Item(468, 6).Value = 100112012
$ws.Cells.Item(468, 7).Value = "Espinaca"
$ws.Cells.Item(468, 8).Value = "Sin especificar"
$ws.Cells.Item(468, 9).Value = "Primera"
$ws.Cells.Item(468, 10).Value = 210
$ws.Cells.Item(468, 11).Value = 5500
$ws.Cells.Item(468, 12).Value = 6000
$ws.Cells.Item(468, 13).Value = 5738
$ws.Cells.Item(468, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(468, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(468, 16).Value = 1913
$ws.Cells.Item(468, 17).Value = 3
$ws.Cells.Item(468, 18).Value = "Hortaliza"
